$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41; this shifts the existing rows 41-76
# down to 42-77 (carrying their formatting/values with them) and updates
# the sheet dimension automatically.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new weekly price entry.
$ws.Range("A41").Value = 11
$ws.Range("B41").Value = "Vega Monumental Concepción"
$ws.Range("C41").Value = "Bíobío"
$ws.Range("D41").Value = 44484
$ws.Range("E41").Value = 8
$ws.Range("F41").Value = 100112043
$ws.Range("G41").Value = "Pepino ensalada"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 450
$ws.Range("K41").Value = 11000
$ws.Range("L41").Value = 12000
$ws.Range("M41").Value = 11556
$ws.Range("N41").Value = "`$/caja 60 unidades"
$ws.Range("O41").Value = "Región de Arica y Parinacota"
$ws.Range("P41").Value = 193
$ws.Range("Q41").Value = 60
$ws.Range("R41").Value = "Hortaliza"
